# Riepilogo_Spese_Edo.xlsx edit
# - Rename the worksheet tab from "Foglio2" to "2025"
# - Re-enter the J2:J13 "Delta spese" formulas (H+I-D_offset pattern) so the
#   block is refreshed/recalculated together, matching the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet tab "Foglio2" -> "2025"
$ws.Name = "2025"

# Re-enter the Delta spese formulas for row 2 through 13 (J2:J13).
# Each cell subtracts the "Costi fissi" value located 15 rows below (D17..D28)
# from the sum of Totale Spese (H) and Spese Previste (I) on the same row.
$ws.Range("J2:J13").FormulaR1C1 = "=RC[-2]+RC[-1]-R[15]C[-6]"

$wb.Application.Calculate()
